# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the "b" file
# has been handed off (new target/handback file + timestamps), and that
# its status moved from "Handed back: in sync with en-US" to
# "Ready for handoff".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": row 3 (b.md) status + datetime columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-42-14 00:42:08"

# ---------------------------------------------------------------------
# Sheet "zh-cn": row 3 (b.md) target file / handback time / status
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-14 00:42:04"

$idx = 0
foreach ($hl in $wsZhCn.Hyperlinks) {
    $idx = $idx + 1
    if ($idx -eq 8) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de": row 3 (b.md) target file / handback time / status
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-14 00:42:08"

$idx = 0
foreach ($hl in $wsDeDe.Hyperlinks) {
    $idx = $idx + 1
    if ($idx -eq 8) {
        $hl.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
